$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the title "Statement of Work" + " - Project 3: Bin Packing" runs
#    into a single run's text (same formatting on both original runs, so a
#    Find/Replace over the combined text collapses them into one <w:r>).
# ---------------------------------------------------------------------------
$enDash = [char]8211
$titleOld = "Statement of Work" + [char]32 + $enDash + " Project 3: Bin Packing"
$d.Content.Find.Execute($titleOld, $true, $false, $false, $false, $false, $true, 1, $false, $titleOld, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Split the "Wrote Bin class" bullet (numId=3) into two bullets:
#      - a new bullet that keeps the original text "Wrote Bin class"
#      - the following bullet becomes "Wrote WinnerTree (Tournament tree)
#        class" with WinnerTree spell-check-flagged, plus the _GoBack
#        bookmark moved to the end of that bullet.
#    The title merge above does not change paragraph count/order, so the
#    "Wrote Bin class" bullet is still paragraph 9.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(9)
$p.Range.InsertParagraphBefore()
$newFirstPara = $d.Paragraphs.Item(9)
$newFirstPara.Range.Text = "Wrote Bin class"

# The original bullet (with its numPr / rPr) is now one paragraph later.
$secondPara = $d.Paragraphs.Item(10)

# Replace just "Bin class" inside that paragraph with the WinnerTree text,
# wrapped in proofErr spell-check markers, via a raw OOXML fragment.
$searchRng = $secondPara.Range.Duplicate()
$searchRng.Find.Execute("Bin class", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$replaceRng = $d.Range($searchRng.Start, $searchRng.End)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$runProps = "<w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>"
$winnerTreeXml = "<w:p $wNs>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r>$runProps<w:t>WinnerTree</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r>$runProps<w:t xml:space='preserve'> (Tournament tree) class</w:t></w:r>" +
    "</w:p>"
$replaceRng.InsertXML($winnerTreeXml)

# Move the _GoBack bookmark to the end of this paragraph's text (collapsed,
# right before the paragraph mark). Adding a bookmark with an existing name
# relocates it, so this also removes it from its old position automatically.
$secondParaAfter = $d.Paragraphs.Item(10)
$secondParaAfter.Range.InsertAfter("X")
$endPos = $secondParaAfter.Range.End - 1
$markerRng = $d.Range($endPos - 1, $endPos)
$d.Bookmarks.Add("_GoBack", $markerRng)
$markerRng.Text = ""

# ---------------------------------------------------------------------------
# 3. Merge "Wrote " + "test cases and testing functions" runs into one run.
# ---------------------------------------------------------------------------
$testText = "Wrote test cases and testing functions"
$d.Content.Find.Execute($testText, $true, $false, $false, $false, $false, $true, 1, $false, $testText, 2) | Out-Null

Write-Output "edit complete"
